# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and target-size parameters were regenerated,
# so every shared string built from them needs the same token swap:
#   D80 -> D86, D64 -> D69, D51 -> D55 (distances)
#   S30 -> S31                        (size)
#
# These tokens only ever show up inside the Condition, Filename_Left,
# Filename_Right, Distance and Size columns (B, D, E, H, J), so walk the
# used range and rewrite any cell whose text contains one of the old
# tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = [ordered]@{
    'D80' = 'D86'
    'D64' = 'D69'
    'D51' = 'D55'
    'S30' = 'S31'
}

$cols = @(2, 4, 5, 8, 10)   # B=Condition, D=Filename_Left, E=Filename_Right, H=Distance, J=Size

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Cells.Item($row, $col)
        $value = $cell.Value2

        if ($value -is [string]) {
            $newValue = $value
            foreach ($old in $map.Keys) {
                $newValue = $newValue.Replace($old, $map[$old])
            }

            if ($newValue -ne $value) {
                $cell.Value = $newValue
            }
        }
    }
}
